$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column B (Encontrados_GitHub) to 0 for the specified rows
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 0
$ws.Range("B15").Value = 0

# Set column C (Encontrados_GitLab) counters for specified rows
$ws.Range("C3").Value = 3
$ws.Range("C5").Value = 2
$ws.Range("C9").Value = 15
$ws.Range("C15").Value = 16
